$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '62.114.68'
$ws.Range("E2").Value = '  +1.51%  '
$ws.Range("D3").Value = '2.416.36'
$ws.Range("E3").Value = '  +1.67%  '
$ws.Range("E4").Value = '  -0.07%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '557.62'
$ws.Range("D5").ClearFormats()
$ws.Range("E5").Value = '  +1.57%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '142.92'
$ws.Range("D6").ClearFormats()
$ws.Range("E6").Value = '  +3.22%  '
$ws.Range("E7").Value = '  +0.07%  '
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.532'
$ws.Range("D8").ClearFormats()
$ws.Range("E8").Value = '  +0.69%  '
$ws.Range("D9").Value = '2.412.12'
$ws.Range("E9").Value = '  +1.44%  '
$ws.Range("E10").Value = '  +0.84%  '
$ws.Range("E11").Value = '  -1.00%  '
$ws.Range("E12").Value = '  +1.37%  '
$ws.Range("E13").Value = '  +1.34%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '26.14'
$ws.Range("D14").ClearFormats()
$ws.Range("E14").Value = '  +4.43%  '
$ws.Range("E15").Value = '  +5.10%  '
$ws.Range("D16").Value = '2.863.65'
$ws.Range("E16").Value = '  +2.68%  '
$ws.Range("D17").Value = '61.929.93'
$ws.Range("E17").Value = '  +1.33%  '
$ws.Range("D18").Value = '2.413.49'
$ws.Range("E18").Value = '  +1.05%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '11.14'
$ws.Range("D19").ClearFormats()
$ws.Range("E19").Value = '  +2.80%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '4.19'
$ws.Range("D20").ClearFormats()
$ws.Range("E20").Value = '  +0.94%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '324.06'
$ws.Range("D21").ClearFormats()
$ws.Range("E21").Value = '  +0.90%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '6.75'
$ws.Range("D22").ClearFormats()
$ws.Range("E22").Value = '  +0.27%  '
$ws.Range("E23").Value = '  +0.00%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '65.37'
$ws.Range("D24").ClearFormats()
$ws.Range("E24").Value = '  +1.53%  '
$ws.Range("E25").Value = '  +1.94%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '9.01'
$ws.Range("D26").ClearFormats()
$ws.Range("E26").Value = '  +7.88%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '586.38'
$ws.Range("D27").ClearFormats()
$ws.Range("E27").Value = '  +15.66%  '
$ws.Range("E28").Value = '  -0.04%  '
$ws.Range("D29").Value = '2.521.49'
$ws.Range("E29").Value = '  +1.38%  '
$ws.Range("D30").Value = '0.0₃0936'
$ws.Range("E30").Value = '  +5.50%  '
$ws.Range("E31").Value = '  +1.47%  '
$ws.Range("E32").Value = '  +5.42%  '
$ws.Range("E33").Value = '  -1.17%  '
$ws.Range("E34").Value = '  +2.34%  '
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '1.56'
$ws.Range("D35").ClearFormats()
$ws.Range("E35").Value = '  +3.07%  '
$ws.Range("E36").Value = '  +5.95%  '
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '4.77'
$ws.Range("D38").ClearFormats()
$ws.Range("E38").Value = '  +2.55%  '
$ws.Range("E39").Value = '  +1.35%  '
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '18.68'
$ws.Range("D40").ClearFormats()
$ws.Range("E40").Value = '  +0.54%  '
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '150.12'
$ws.Range("D41").ClearFormats()
$ws.Range("E41").Value = '  +2.28%  '
$ws.Range("E42").Value = '  -2.41%  '
$ws.Range("E43").Value = '  +0.03%  '
$ws.Range("E44").Value = '  +12.28%  '
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '151.13'
$ws.Range("D45").ClearFormats()
$ws.Range("E45").Value = '  +2.00%  '
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '3.66'
$ws.Range("D46").ClearFormats()
$ws.Range("E46").Value = '  +1.57%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '0.0540'
$ws.Range("D47").ClearFormats()
$ws.Range("E47").Value = '  +3.86%  '
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '20.22'
$ws.Range("D48").ClearFormats()
$ws.Range("E49").Value = '  +2.45%  '
$ws.Range("E50").Value = '  +1.35%  '
$ws.Range("E51").Value = '  +2.01%  '
